$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.915.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.68%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.350.96"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.04%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.674"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -3.01%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'240.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.59%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'72.88"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -4.92%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.03%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.599"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.14%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -2.93%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'58.88"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.52%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'33.33"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +2.99%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +0.05%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'7.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -3.00%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'2.699.62"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.17%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'16.35"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -5.11%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -2.09%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.345.18"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.39%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'43.816.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -1.03%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.0000103"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.67%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'6.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.37%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'78.33"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.62%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'255.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.37%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'1.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +9.09%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -0.03%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D27").Value = "'2.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -2.75%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'10.54"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -3.14%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'2.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.73%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'22.56"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -2.55%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'177.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +1.05%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -1.83%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +0.95%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.0751"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.89%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("B35").Value = "'InternetComputer(DFINITY)"
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = "'5.46"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +2.14%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("B36").Value = "'Filecoin"
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = "'5.14"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -5.17%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'3.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -3.46%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'6.45"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -2.28%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -4.76%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -0.54%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'67.81"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +25.35%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'5.11"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +13.92%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +8.04%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +1.17%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +3.56%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'18.92"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.78%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'2.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.52%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -1.83%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +0.15%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'99.42"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -2.98%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -5.20%  "
$ws.Range("E51").Style = "Normal"
